$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 3 ("SMB_loading"), pushing existing rows down.
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 "JAN t+1"
$ws.Cells.Item(3, 1).Value = "JAN t+1"
$ws.Cells.Item(3, 2).Value = 0.01826231019968784
$ws.Cells.Item(3, 3).Value = 0.0155582838208262
$ws.Cells.Item(3, 4).Value = 1.173799784732172
$ws.Cells.Item(3, 5).Value = 0.2415370822754689

# Copy style from row 2 (A column style s="1") to the new row's A cell
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(3, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(3, 1).Value = "JAN t+1"

# Update all other values per diff (rows shifted down by one starting from old row 3)

# Row 2: intercept
$ws.Cells.Item(2, 2).Value = 0.02963681848556416
$ws.Cells.Item(2, 3).Value = 0.01814905270938225
$ws.Cells.Item(2, 4).Value = 1.632967789566408
$ws.Cells.Item(2, 5).Value = 0.1036723649356633

# Row 4: SMB_loading
$ws.Cells.Item(4, 2).Value = 0.003171075620851749
$ws.Cells.Item(4, 3).Value = 0.001615447416179446
$ws.Cells.Item(4, 4).Value = 1.962970499127346
$ws.Cells.Item(4, 5).Value = 0.05070345951672352

# Row 5: HML_loading
$ws.Cells.Item(5, 2).Value = 0.0006577907896868769
$ws.Cells.Item(5, 3).Value = 0.0009707244012964668
$ws.Cells.Item(5, 4).Value = 0.6776287778573956
$ws.Cells.Item(5, 5).Value = 0.4986027694196633

# Row 6: RMW_loading
$ws.Cells.Item(6, 2).Value = -0.00009124669264801628
$ws.Cells.Item(6, 3).Value = 0.001682365601839705
$ws.Cells.Item(6, 4).Value = -0.05423713641567324
$ws.Cells.Item(6, 5).Value = 0.9567874101992

# Row 7: CMA_loading
$ws.Cells.Item(7, 2).Value = -0.0005290258431414854
$ws.Cells.Item(7, 3).Value = 0.001009039239031843
$ws.Cells.Item(7, 4).Value = -0.5242866904255152
$ws.Cells.Item(7, 5).Value = 0.6005207337994569

# Row 8: mkt_loading
$ws.Cells.Item(8, 2).Value = -0.0003026426544363091
$ws.Cells.Item(8, 3).Value = 0.001416311050514604
$ws.Cells.Item(8, 4).Value = -0.2136837485849925
$ws.Cells.Item(8, 5).Value = 0.8309592460526546

# Row 9: size*JAN_loading
$ws.Cells.Item(9, 2).Value = -0.007307766519562536
$ws.Cells.Item(9, 3).Value = 0.54696310770923
$ws.Cells.Item(9, 4).Value = -0.01336062051820761
$ws.Cells.Item(9, 5).Value = 0.9893502136755956

# Row 10: BM*JAN_loading
$ws.Cells.Item(10, 2).Value = 0.538102921313061
$ws.Cells.Item(10, 3).Value = 0.4995949554749016
$ws.Cells.Item(10, 4).Value = 1.0770783720215
$ws.Cells.Item(10, 5).Value = 0.2824322186534419

# Row 11: ROE*JAN_loading
$ws.Cells.Item(11, 2).Value = -0.0001979432691203469
$ws.Cells.Item(11, 3).Value = 0.002290959030434774
$ws.Cells.Item(11, 4).Value = -0.0864019244738661
$ws.Cells.Item(11, 5).Value = 0.931212659957

# Row 12: INV*JAN_loading (new row, was not present before at this position)
$ws.Cells.Item(12, 1).Value = "INV*JAN_loading"
$ws.Cells.Item(12, 2).Value = 0.0009285788881154748
$ws.Cells.Item(12, 3).Value = 0.002120422863847711
$ws.Cells.Item(12, 4).Value = 0.4379215598677705
$ws.Cells.Item(12, 5).Value = 0.6618025878139184

# Row 13: mkt_risk_premium*JAN_loading
$ws.Cells.Item(13, 1).Value = "mkt_risk_premium*JAN_loading"
$ws.Cells.Item(13, 2).Value = -0.2109742443719655
$ws.Cells.Item(13, 3).Value = 0.296482670251288
$ws.Cells.Item(13, 4).Value = -0.7115904757372542
$ws.Cells.Item(13, 5).Value = 0.4773490894200991
